$d = $word.ActiveDocument

$replacements = @(
    @{Old = "31×48="; New = "12×26="},
    @{Old = "46×39="; New = "27×17="},
    @{Old = "36×57="; New = "73×31="},
    @{Old = "79×18="; New = "35×33="},
    @{Old = "48×16="; New = "79×46="},
    @{Old = "80×90="; New = "61×14="},
    @{Old = "52×43="; New = "51×70="},
    @{Old = "85×57="; New = "52×44="},
    @{Old = "31×84="; New = "41×81="},
    @{Old = "92×63="; New = "16×40="},
    @{Old = "39×30="; New = "51×55="},
    @{Old = "21×49="; New = "96×99="},
    @{Old = "16×13="; New = "63×40="},
    @{Old = "35×17="; New = "69×95="},
    @{Old = "15×92="; New = "26×54="},
    @{Old = "68×24="; New = "43×44="},
    @{Old = "27×54="; New = "19×90="},
    @{Old = "46×93="; New = "96×68="},
    @{Old = "74×89="; New = "52×88="},
    @{Old = "98×47="; New = "67×42="},
    @{Old = "87×33="; New = "12×93="},
    @{Old = "13×62="; New = "16×75="},
    @{Old = "81×78="; New = "57×85="},
    @{Old = "30×88="; New = "11×39="},
    @{Old = "95×71="; New = "43×32="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.Old, $true, $true, $false, $false, $false, $true, 1, $false, $r.New, 2)
}
